# Apply the "Updated cryptos list" refresh (prices / 1h-volume deltas, plus a
# row-order swap between InjectiveProtocol and Mantle) described by the diff.
#
# All target cells in columns D/E (and the B/C text swap in rows 45-46) are
# stored as literal text in the workbook (t="inlineStr"). Column D sometimes
# holds values that *look* numeric (e.g. "0.488", "2.31"); a plain
# `Range.Value = "0.488"` assignment lets Excel's COM layer auto-coerce that
# into a real floating-point number (and can introduce binary-float noise,
# e.g. 2.3100000000000001). To keep those cells as exact text - matching the
# source workbook - we prefix them with a leading apostrophe (the standard
# "force text" entry trick), then reset the cell Style back to "Normal" so we
# don't leave a quotePrefix style applied that the original file never had.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($addr, $value)
    $ws.Range($addr).Value = $value
}

function Set-NumericLookingText {
    param($addr, $value)
    # Leading apostrophe forces Excel to keep the entry as text instead of
    # parsing it into a number.
    $ws.Range($addr).Value = "'" + $value
}

# --- Row 2 (Bitcoin) ---
Set-TextValue "D2" "67.582.39"
Set-TextValue "E2" "  +0.27%  "

# --- Row 3 (Ethereum) ---
Set-TextValue "D3" "3.505.23"
Set-TextValue "E3" "  -0.45%  "

# --- Row 4 (TetherUSD) ---
Set-TextValue "E4" "  -0.01%  "

# --- Row 5 (BNB) ---
Set-NumericLookingText "D5" "605.93"
Set-TextValue "E5" "  -1.20%  "

# --- Row 6 (Solana) ---
Set-NumericLookingText "D6" "151.85"
Set-TextValue "E6" "  +0.36%  "

# --- Row 7 (LidoStakedEther) ---
Set-TextValue "D7" "3.504.41"
Set-TextValue "E7" "  -0.44%  "

# --- Row 8 (USDC) ---
Set-TextValue "E8" "  -0.03%  "

# --- Row 9 (XRP) ---
Set-NumericLookingText "D9" "0.488"
Set-TextValue "E9" "  +1.37%  "

# --- Row 10 (Dogecoin) ---
Set-NumericLookingText "D10" "0.144"
Set-TextValue "E10" "  +2.65%  "

# --- Row 11 (Toncoin) ---
Set-NumericLookingText "D11" "7.63"
Set-TextValue "E11" "  +6.93%  "

# --- Row 12 (Cardano) ---
Set-NumericLookingText "D12" "0.432"
Set-TextValue "E12" "  +1.69%  "

# --- Row 13 (ShibaInu) ---
Set-TextValue "E13" "  -1.96%  "

# --- Row 14 (Avalanche) ---
Set-NumericLookingText "D14" "32.17"
Set-TextValue "E14" "  +0.37%  "

# --- Row 15 (WrappedliquidstakedEther2.0) ---
Set-TextValue "D15" "4.097.67"
Set-TextValue "E15" "  -0.48%  "

# --- Row 16 (WrappedEther) ---
Set-TextValue "D16" "3.506.36"
Set-TextValue "E16" "  -0.47%  "

# --- Row 17 (WrappedBTC) ---
Set-TextValue "D17" "67.515.29"
Set-TextValue "E17" "  +0.17%  "

# --- Row 18 (TRON) ---
Set-TextValue "E18" "  -0.57%  "

# --- Row 19 (Polkadot) ---
Set-NumericLookingText "D19" "6.49"
Set-TextValue "E19" "  +1.70%  "

# --- Row 20 (Chainlink) ---
Set-NumericLookingText "D20" "15.46"
Set-TextValue "E20" "  +1.10%  "

# --- Row 21 (Uniswap) ---
Set-NumericLookingText "D21" "9.84"
Set-TextValue "E21" "  +2.82%  "

# --- Row 22 (BitcoinCash) ---
Set-NumericLookingText "D22" "447.29"
Set-TextValue "E22" "  +0.40%  "

# --- Row 23 (Polygon) ---
Set-TextValue "E23" "  +0.82%  "

# --- Row 24 (Litecoin) ---
Set-NumericLookingText "D24" "78.31"
Set-TextValue "E24" "  +1.13%  "

# --- Row 25 (WrappedeETH) ---
Set-TextValue "D25" "3.647.55"
Set-TextValue "E25" "  -0.41%  "

# --- Row 26 (PEPE) ---
Set-TextValue "E26" "  -3.73%  "

# --- Row 28 (RenderToken) ---
Set-TextValue "E28" "  +3.16%  "

# --- Row 29 (InternetComputer(DFINITY)) ---
Set-NumericLookingText "D29" "10.03"
Set-TextValue "E29" "  -1.50%  "

# --- Row 30 (PancakeSwap) ---
Set-NumericLookingText "D30" "2.51"
Set-TextValue "E30" "  -0.14%  "

# --- Row 31 (Fetch.AI) ---
Set-TextValue "E31" "  +4.95%  "

# --- Row 32 (Kaspa) ---
Set-NumericLookingText "D32" "0.173"
Set-TextValue "E32" "  +5.11%  "

# --- Row 33 (Binance-PegBSC-USD) ---
Set-TextValue "E33" "  +0.00%  "

# --- Row 34 (EthereumClassic) ---
Set-TextValue "E34" "  -0.87%  "

# --- Row 35 (NEARProtocol) ---
Set-NumericLookingText "D35" "6.15"
Set-TextValue "E35" "  -0.27%  "

# --- Row 37 (RenzoRestakedETH) ---
Set-TextValue "D37" "3.495.99"

# --- Row 38 (Aptos) ---
Set-TextValue "E38" "  -0.24%  "

# --- Row 40 (Stacks) ---
Set-NumericLookingText "D40" "2.31"
Set-TextValue "E40" "  +6.86%  "

# --- Row 41 (Monero) ---
Set-NumericLookingText "D41" "179.74"
Set-TextValue "E41" "  +2.48%  "

# --- Row 42 (FirstDigitalUSD) ---
Set-NumericLookingText "D42" "0.999"
Set-TextValue "E42" "  -0.04%  "

# --- Row 43 (Hedera) ---
Set-TextValue "E43" "  +1.51%  "

# --- Row 44 (Filecoin) ---
Set-TextValue "E44" "  +0.59%  "

# --- Row 45 / 46 swap: InjectiveProtocol now ranks above Mantle ---
Set-TextValue "B45" "InjectiveProtocol"
Set-TextValue "C45" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-NumericLookingText "D45" "30.32"
Set-TextValue "E45" "  +7.38%  "

Set-TextValue "B46" "Mantle"
Set-TextValue "C46" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-NumericLookingText "D46" "0.893"
Set-TextValue "E46" "  +1.48%  "

# --- Row 47 (OKB) ---
Set-NumericLookingText "D47" "46.48"
Set-TextValue "E47" "  +3.11%  "

# --- Row 48 (ONDO) ---
Set-TextValue "E48" "  +4.33%  "

# --- Row 49 (dogwifhat) ---
Set-NumericLookingText "D49" "2.55"
Set-TextValue "E49" "  -3.51%  "

# --- Row 50 (Cosmos) ---
Set-NumericLookingText "D50" "7.62"
Set-TextValue "E50" "  +0.30%  "

# --- Row 51 (TheGraph) ---
Set-TextValue "E51" "  +1.73%  "

# Clean up the quote-prefix style that forcing-to-text left behind, so the
# affected cells keep the workbook's default (unstyled) look, same as before.
$textForcedCells = @(
    "D5","D6","D9","D10","D11","D12","D14","D19","D20","D21","D22","D24",
    "D29","D30","D32","D35","D40","D41","D42","D45","D46","D47","D49","D50"
)
foreach ($addr in $textForcedCells) {
    $ws.Range($addr).Style = "Normal"
}
